$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "problem 10" solution statistics (new data row added/updated) ---
$ws.Range("B2").Value = 15867
$ws.Range("C2").Value = 4974

$ws.Range("B3").Value = 11637
$ws.Range("C3").Value = 1636

$ws.Range("B4").Value = 9386
$ws.Range("C4").Value = 947

$ws.Range("B5").Value = 8231
$ws.Range("C5").Value = 418

$ws.Range("B6").Value = 6558
$ws.Range("C6").Value = 1268

$ws.Range("B7").Value = 5498
$ws.Range("C7").Value = 247

$ws.Range("B8").Value = 3372
$ws.Range("C8").Value = 128

$ws.Range("B9").Value = 3009
$ws.Range("C9").Value = 177

$ws.Range("B10").Value = 2110
$ws.Range("C10").Value = 55

# E10 used to be a placeholder "-" (no prior day to compare against); now that
# row 9 -> row 10 (day 9 -> day 10) has real data, give it the same "Var"
# percentage-change formula/format used by the rows above it.
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Formula = "=-(1-(D10/D9))"
$excel.CutCopyMode = 0

# --- Nudge the chart slightly to the right/down as in the original edit ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 295.1874803149607
$co.Top = 0.75
$co.Width = 803.9374803149606
$co.Height = 384.7500787401575

# --- Add the new (empty) Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet2"

# --- Restore selection on Sheet1 and make it the active sheet ---
$ws.Select()
$ws.Range("U6").Select()
